# Generate Report for handback
# - Flip status text from "Ready for handoff" to "Handed back: in sync with en-US"
#   (applies everywhere that text is used: Overview sheet + the per-locale sheets)
# - For the zh-cn and de-de sheets, record the handback results for the two
#   in-scope rows (row 2 = a.md.md, row 3 = b.md.md): populate the
#   "Latest Target File" (E) / "Latest Handback File" (F) columns with
#   hyperlinked file names, and stamp the "Latest Handback DateTime" (G)
#   column with the handback timestamp.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Update the status text everywhere it appears ("Ready for handoff" -> "Handed back: in sync with en-US")
# ---------------------------------------------------------------------------
foreach ($s in $wb.Worksheets) {
    $used = $s.UsedRange
    foreach ($cell in $used.Cells) {
        if ($cell.Value() -eq "Ready for handoff") {
            $cell.Value = "Handed back: in sync with en-US"
        }
    }
}

# ---------------------------------------------------------------------------
# 2. zh-cn sheet
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$zhSourceUrl = "https://github.com/OpenLocalizationTest/oltest/blob/52f3753b6a36c63a252b14a1bd1c568084d5fee2/e2e/a.md.md"
$zhHandoffUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/a3f5ce1547834005309ca6f1fa27bcbeef505136/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/a.md.370104d57010292b5765347db07350cde3a977e8.zh-cn.xlf"
$zhHandoffFileName = "a.md.370104d57010292b5765347db07350cde3a977e8.zh-cn.xlf"

$wsZh.Hyperlinks.Add($wsZh.Range("E2"), $zhSourceUrl, "", "", "a.md.md")
$wsZh.Hyperlinks.Add($wsZh.Range("F2"), $zhHandoffUrl, "", "", $zhHandoffFileName)
$wsZh.Hyperlinks.Add($wsZh.Range("E3"), $zhSourceUrl, "", "", "a.md.md")
$wsZh.Hyperlinks.Add($wsZh.Range("F3"), $zhHandoffUrl, "", "", $zhHandoffFileName)

$wsZh.Range("G2").Value = "2016-02-16 14:49:16"
$wsZh.Range("G3").Value = "2016-02-16 14:49:16"

# ---------------------------------------------------------------------------
# 3. de-de sheet
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$deSourceUrl = "https://github.com/OpenLocalizationTest/oltest/blob/52f3753b6a36c63a252b14a1bd1c568084d5fee2/e2e/a.md.md"
$deHandoffUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/f1274742262761b61fd2cbcf061a3e19f3505b8a/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/a.md.370104d57010292b5765347db07350cde3a977e8.de-de.xlf"
$deHandoffFileName = "a.md.370104d57010292b5765347db07350cde3a977e8.de-de.xlf"

$wsDe.Hyperlinks.Add($wsDe.Range("E2"), $deSourceUrl, "", "", "a.md.md")
$wsDe.Hyperlinks.Add($wsDe.Range("F2"), $deHandoffUrl, "", "", $deHandoffFileName)
$wsDe.Hyperlinks.Add($wsDe.Range("E3"), $deSourceUrl, "", "", "a.md.md")
$wsDe.Hyperlinks.Add($wsDe.Range("F3"), $deHandoffUrl, "", "", $deHandoffFileName)

$wsDe.Range("G2").Value = "2016-02-16 14:49:44"
$wsDe.Range("G3").Value = "2016-02-16 14:49:44"

Write-Output "Handback report generated"
